# Update the arithmetic problems in the three-digit x one-digit multiplication
# worksheet to reflect the newly generated values.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "658×6="; New = "852×2=" },
    @{ Old = "799×6="; New = "544×3=" },
    @{ Old = "850×2="; New = "225×2=" },
    @{ Old = "398×9="; New = "622×8=" },
    @{ Old = "367×6="; New = "838×8=" },
    @{ Old = "867×4="; New = "219×9=" },
    @{ Old = "261×2="; New = "448×6=" },
    @{ Old = "749×4="; New = "364×7=" },
    @{ Old = "380×2="; New = "399×7=" },
    @{ Old = "326×4="; New = "904×2=" },
    @{ Old = "817×7="; New = "669×6=" },
    @{ Old = "301×3="; New = "128×4=" },
    @{ Old = "528×5="; New = "169×7=" },
    @{ Old = "460×8="; New = "143×4=" },
    @{ Old = "837×3="; New = "620×7=" },
    @{ Old = "429×2="; New = "619×6=" },
    @{ Old = "893×3="; New = "305×7=" },
    @{ Old = "741×7="; New = "321×9=" },
    @{ Old = "269×3="; New = "857×9=" },
    @{ Old = "744×2="; New = "516×2=" },
    @{ Old = "335×3="; New = "522×7=" },
    @{ Old = "525×2="; New = "265×8=" },
    @{ Old = "490×4="; New = "182×6=" },
    @{ Old = "269×7="; New = "182×2=" },
    @{ Old = "223×8="; New = "878×9=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
